# Fill in the missing "Reference" values ("-") for the last three rows of
# the "Tabela1" table (ONI, PDO, QBO), which previously had an empty
# Reference (column G) cell, and move the active selection to G27 (the
# cell just below the table) to reflect where the user ended up after
# typing the values in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G24").Value = "-"
$ws.Range("G25").Value = "-"
$ws.Range("G26").Value = "-"

$ws.Range("G27").Select()
